$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.593.26'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '3.509.89'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''609.48'
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").Value = '''152.05'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = '3.509.09'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = '''0.487'
$ws.Range("E9").Value = '  +0.93%  '
$ws.Range("E10").Value = '  +2.66%  '
$ws.Range("D11").Value = '''7.66'
$ws.Range("E11").Value = '  +7.37%  '
$ws.Range("D12").Value = '''0.433'
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = '''32.22'
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '4.101.43'
$ws.Range("D16").Value = '3.510.86'
$ws.Range("E16").Value = '  -0.74%  '
$ws.Range("D17").Value = '67.559.23'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = '''6.51'
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("D21").Value = '''9.87'
$ws.Range("E21").Value = '  +3.71%  '
$ws.Range("D22").Value = '''447.38'
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("E23").Value = '  +0.73%  '
$ws.Range("D24").Value = '''78.34'
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("D25").Value = '3.649.53'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").Value = '''0.0000128'
$ws.Range("E26").Value = '  -3.45%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = '''8.78'
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("D29").Value = '''10.05'
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("D30").Value = '''2.52'
$ws.Range("E30").Value = '  +0.35%  '
$ws.Range("E31").Value = '  +4.83%  '
$ws.Range("D32").Value = '''0.173'
$ws.Range("E32").Value = '  +5.74%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = '''25.62'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("E36").Value = '  +0.86%  '
$ws.Range("D37").Value = '3.499.26'
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("D38").Value = '''8.00'
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("E40").Value = '  +7.18%  '
$ws.Range("D41").Value = '''179.23'
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +1.67%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '''30.43'
$ws.Range("E45").Value = '  +7.96%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.895'
$ws.Range("E46").Value = '  +1.53%  '
$ws.Range("D47").Value = '''46.44'
$ws.Range("E47").Value = '  +2.96%  '
$ws.Range("D48").Value = '''1.30'
$ws.Range("E48").Value = '  +4.11%  '
$ws.Range("D49").Value = '''2.55'
$ws.Range("E49").Value = '  -3.32%  '
$ws.Range("D50").Value = '''7.62'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("E51").Value = '  +1.76%  '
